# Plan de Proyecto - Kairos - NexTech
# Applies the changes described by the commit:
#   docs(Etapa Inicio): actualizar documentos
#
# 1) Rename four floating pictures in the document body (Shape.Name mirrors
#    the wp:docPr/@name the Word object model exposes for anchored shapes).
# 2) Rename two floating pictures that live in the page header.
# 3) Tweak the risk-management paragraph wording.
#
# (The ToC content control's internal w:id is a Word-managed, read-only
# value - ContentControl.ID has no setter in the Word object model - so it
# is intentionally left alone here.)

$d = $word.ActiveDocument

# --- Body pictures --------------------------------------------------------
for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    $shp = $d.Shapes.Item($i)

    if ($shp.Name -eq "image5.png") {
        $shp.Name = "image2.png"
    }
    elseif ($shp.Name -eq "image3.jpg") {
        $shp.Name = "image5.jpg"
    }
    elseif ($shp.Name -eq "image4.png") {
        $shp.Name = "image1.png"
    }
}

# --- Header pictures -------------------------------------------------------
$hdr = $d.Sections.Item(1).Headers.Item(1)
for ($i = 1; $i -le $hdr.Shapes.Count; $i++) {
    $shp = $hdr.Shapes.Item($i)

    if ($shp.Name -eq "image2.png") {
        $shp.Name = "image3.png"
    }
    elseif ($shp.Name -eq "image1.jpg") {
        $shp.Name = "image4.jpg"
    }
}

# --- Risk-management paragraph wording -------------------------------------
$d.Content.Find.Execute(
    "el sistema para gestión de los riesgos Vesta Risk.", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "el sistema para gestión de riesgos, Vesta Risk Manager.", 2) | Out-Null
